$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Data edit: Study ID S005 (row 6) age changed from 27 to 29.
$ws.Range("B6").Value = 29

# C6/D6 hold shared formulas referencing B6 - they recalc automatically.

# Reflect the new selection left by the editor (was G10, now B6).
$ws.Range("B6").Select()
